# Generate Report for Handback
# Update the timestamp values recorded on the handback status report.
$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-31 01:10:13"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for the first row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-31 01:09:59"
$wsZhCn.Range("K2").Value = "2016-08-31 01:10:32"

# de-de sheet: "Correspond Handback DateTime" for the first row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-31 01:10:39"
